# Horarios Línea 141 - actualización 08:11:18
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 08:11:18"
$ws1.Range("A3").Value = "Total filas: 62"

# Rows 20 & 21 swap their Hora_Scrap / Linea / Minutos values
$ws1.Range("A20").Value = "05:57:04"
$ws1.Range("C20").Value = "16_SANTA ANA"
$ws1.Range("D20").Value = 84

$ws1.Range("A21").Value = "06:16:41"
$ws1.Range("C21").Value = "23_HERNANDEZ"
$ws1.Range("D21").Value = 65

# Rows 50-61 get reshuffled (new scrap cycle interleaved with previous data)
# row, Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada
$rows1 = @(
    @(50, "08:11:18", "09:05", "10_OLMOS",            54,  "LP1912"),
    @(51, "07:38:09", "09:06", "16_SANTA ANA",        88,  "LP1912"),
    @(52, "07:56:02", "09:11", "27_EL RETIRO",        75,  "LP1912"),
    @(53, "07:49:32", "09:12", "27_EL RETIRO",        83,  "LP1912"),
    @(54, "07:38:09", "09:14", "27_EL RETIRO",        96,  "LP1912"),
    @(55, "07:49:32", "09:17", "14_ABASTO",           88,  "LP1912"),
    @(56, "07:38:09", "09:18", "15X38_ABASTO",        100, "LP1912"),
    @(57, "07:38:09", "09:18", "14_ABASTO",           100, "LP1912"),
    @(58, "08:11:18", "09:28", "23_HERNANDEZ",        77,  "LP1912"),
    @(59, "07:38:09", "09:29", "10_OLMOS",            111, "LP1912"),
    @(60, "08:11:18", "09:31", "16_SANTA ANA",        80,  "LP1912"),
    @(61, "07:49:32", "09:39", "15_ABASTO",           110, "LP1912"),
    @(62, "07:49:32", "09:41", "11_ETCHEVERRY",       112, "LP1912"),
    @(63, "07:56:02", "09:42", "11_ETCHEVERRY",       106, "LP1912"),
    @(64, "07:49:32", "09:43", "16_P MOR-SANTA ANA",  114, "LP1912"),
    @(65, "08:11:18", "09:53", "10_OLMOS",            102, "LP1912"),
    @(66, "08:11:18", "09:59", "215C_EL PATO",        108, "LP1912"),
    @(67, "08:11:18", "10:06", "14_ABASTO",           115, "LP1912")
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Range("A$r").Value = $row[1]
    $ws1.Range("B$r").Value = $row[2]
    $ws1.Range("C$r").Value = $row[3]
    $ws1.Range("D$r").Value = $row[4]
    $ws1.Range("E$r").Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 08:11:18"
$ws2.Range("A3").Value = "Total filas: 14"

$ws2.Range("A19").Value = "08:11:18"
$ws2.Range("B19").Value = "09:59"
$ws2.Range("C19").Value = "215C_EL PATO"
$ws2.Range("D19").Value = 108
$ws2.Range("E19").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 08:11:18"
